# Updates cryptos list (prices / 1h volume change) and reorders a few
# coin rows, matching the "Updated cryptos list ... with GitHub Actions"
# commit. Column D (Price) values are entered with a leading apostrophe
# so Excel keeps them as literal text (e.g. "1.00", "95.815.20") instead
# of silently re-parsing them as numbers and dropping formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.815.20"
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = "'3.568.52"
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'236.15"
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").Value = "'651.85"
$ws.Range("E6").Value = '  +3.50%  '
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'0.998"
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("D11").Value = "'3.568.62"
$ws.Range("E11").Value = '  +1.66%  '
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = "'42.31"
$ws.Range("E13").Value = '  -3.10%  '
$ws.Range("D14").Value = "'6.53"
$ws.Range("E14").Value = '  +3.91%  '
$ws.Range("D15").Value = "'4.256.61"
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("D16").Value = "'95.657.20"
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = "'3.572.24"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("E19").Value = '  -5.41%  '
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("D21").Value = "'17.86"
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("E22").Value = '  +3.25%  '
$ws.Range("D23").Value = "'507.65"
$ws.Range("E23").Value = '  -2.31%  '
$ws.Range("D24").Value = "'0.480"
$ws.Range("E24").Value = '  -4.08%  '
$ws.Range("E25").Value = '  +3.69%  '
$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = '  -2.29%  '
$ws.Range("D27").Value = "'95.70"
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").Value = "'12.53"
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("D29").Value = "'3.761.04"
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("D30").Value = "'3.07"
$ws.Range("E30").Value = '  +4.20%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = "'0.140"
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = "'11.25"
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = '  -2.19%  '
$ws.Range("D36").Value = "'32.06"
$ws.Range("E36").Value = '  +7.41%  '
$ws.Range("D37").Value = "'0.561"
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = "'8.20"
$ws.Range("E38").Value = '  +8.00%  '
$ws.Range("D39").Value = "'560.37"
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = "'1.77"
$ws.Range("E44").Value = '  +3.33%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'35.12"
$ws.Range("E45").Value = '  +38.28%  '
$ws.Range("D46").Value = "'2.30"
$ws.Range("E46").Value = '  +5.71%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = "'5.68"
$ws.Range("E47").Value = '  +2.16%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Value = "'23.61"
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("D49").Value = "'0.0413"
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("D51").Value = "'53.56"
$ws.Range("E51").Value = '  -0.70%  '
